$wb = $excel.ActiveWorkbook

# "optimization_parameters" sheet had a stray leftover row (label "Sheet", values 3/4)
# sitting between the expression_timepoints/Strain rows and the simulation_timepoints
# row. Remove it -- everything below shifts up by one row.
$wsParams = $wb.Worksheets.Item("optimization_parameters")
$wsParams.Rows.Item(16).Delete()

# Select the (now shifted-up) simulation_timepoints row so the saved view matches.
$wsParams.Rows.Item(16).Select() | Out-Null

# Move the active tab / selected sheet from optimization_parameters to threshold_b.
$wsThreshold = $wb.Worksheets.Item("threshold_b")
$wsThreshold.Activate()
$wsThreshold.Range("A2").Select() | Out-Null
